$d = $word.ActiveDocument

# The edit removes two whole paragraphs near the end of the "Requisitos"
# section: the "Ver no Jupiter..." paragraph and the "(c) 2020..." footer
# paragraph, while leaving the paragraphs around them (including the
# trailing empty paragraph and the page-break paragraph) untouched.

function Find-ParagraphByExactText($doc, $searchText) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd()
        if ($t -eq $searchText) {
            return $p
        }
    }
    return $null
}

$targets = @(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    [char]0xA9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
)

# Resolve the paragraphs first (by content), then delete them starting
# from the one that occurs latest in the document so earlier ranges stay
# valid.
$paragraphsToDelete = @()
foreach ($text in $targets) {
    $p = Find-ParagraphByExactText $d $text
    if ($p -ne $null) {
        $paragraphsToDelete += , @{ Start = $p.Range.Start; End = $p.Range.End }
    }
}

$paragraphsToDelete = $paragraphsToDelete | Sort-Object { $_.Start } -Descending

foreach ($info in $paragraphsToDelete) {
    $r = $d.Range($info.Start, $info.End)
    $r.Delete()
}
